$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 2.3
$ws.Range("I2").Value = 3.7
$ws.Range("J2").Value = 3.2
$ws.Range("L2").Value = 4.5
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 1.62
$ws.Range("P2").Value = 2.2
$ws.Range("Q2").Value = 3.1
$ws.Range("R2").Value = 1.36
$ws.Range("S2").Value = 1.67
$ws.Range("T2").Value = 2.1
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.53
$ws.Range("W2").Value = 5.5
$ws.Range("X2").Value = 9
$ws.Range("Z2").Value = 21
$ws.Range("AC2").Value = 5
$ws.Range("AE2").Value = 21
$ws.Range("AH2").Value = 17
$ws.Range("AI2").Value = 15
$ws.Range("AK2").Value = 41
$ws.Range("AN2").Value = 4
$ws.Range("AT2").Value = 2.1
$ws.Range("AU2").Value = 10
$ws.Range("AX2").Value = 23
$ws.Range("AY2").Value = 41
$ws.Range("AZ2").Value = 81
$ws.Range("BA2").Value = 151
$ws.Range("G3").Value = 2.8
$ws.Range("I3").Value = 2.88
$ws.Range("J3").Value = 3.6
$ws.Range("L3").Value = 3.6
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("Y3").Value = 12
$ws.Range("Z3").Value = 29
$ws.Range("AA3").Value = 29
$ws.Range("AE3").Value = 19
$ws.Range("AH3").Value = 12
$ws.Range("AP3").Value = 34
$ws.Range("AV3").Value = 81
$ws.Range("AW3").Value = 4.5
$ws.Range("J7").Value = 2.87
$ws.Range("L7").Value = 3.8
$ws.Range("U7").Value = 1.85
$ws.Range("V7").Value = 1.75
$ws.Range("W7").Value = 6.3
$ws.Range("Y7").Value = 9.25
$ws.Range("AI7").Value = 11.5
$ws.Range("AL7").Value = 40
$ws.Range("AU7").Value = 6.8
$ws.Range("AY7").Value = 24
$ws.Range("AZ7").Value = 90
$ws.Range("BB7").Value = 300
$ws.Range("G8").Value = 2.35
$ws.Range("I8").Value = 2.9
$ws.Range("S8").Value = 1.5
$ws.Range("T8").Value = 2.5
$ws.Range("AH8").Value = 15
$ws.Range("AL8").Value = 41
$ws.Range("AN8").Value = 4.33
$ws.Range("AT8").Value = 2.5
$ws.Range("O12").Value = 1.3
$ws.Range("P12").Value = 3.4
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 1.85
$ws.Range("G13").Value = 1.65
$ws.Range("H13").Value = 3.35
$ws.Range("I13").Value = 5.4
$ws.Range("J13").Value = 2.18
$ws.Range("K13").Value = 2.07
$ws.Range("L13").Value = 5.6
$ws.Range("M13").Value = 1.02
$ws.Range("N13").Value = 7.1
$ws.Range("P13").Value = 2.6
$ws.Range("Q13").Value = 2.07
$ws.Range("R13").Value = 1.6
$ws.Range("S13").Value = 1.42
$ws.Range("T13").Value = 2.45
$ws.Range("U13").Value = 2.02
$ws.Range("V13").Value = 1.62
$ws.Range("X13").Value = 6.9
$ws.Range("Z13").Value = 12.5
$ws.Range("AA13").Value = 14.5
$ws.Range("AC13").Value = 7.6
$ws.Range("AD13").Value = 6.7
$ws.Range("AE13").Value = 19.5
$ws.Range("AF13").Value = 120
$ws.Range("AH13").Value = 30
$ws.Range("AI13").Value = 18.5
$ws.Range("AJ13").Value = 120
$ws.Range("AK13").Value = 70
$ws.Range("AL13").Value = 80
$ws.Range("AN13").Value = 3.3
$ws.Range("AO13").Value = 7.9
$ws.Range("AP13").Value = 18.5
$ws.Range("AQ13").Value = 26
$ws.Range("AR13").Value = 60
$ws.Range("AS13").Value = 250
$ws.Range("AU13").Value = 7.9
$ws.Range("AV13").Value = 90
$ws.Range("AW13").Value = 6.8
$ws.Range("AX13").Value = 35
$ws.Range("AY13").Value = 40
$ws.Range("AZ13").Value = 250
$ws.Range("BA13").Value = 300
